$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity Log")
if ($null -eq $ws) { $ws = $wb.ActiveSheet }

# E51: add endtime value
$ws.Range("E51").Value = 0.64930555555555558

# Row 52
$ws.Range("B52").Value = 6977
$ws.Range("C52").Value = 43926
$ws.Range("D52").Value = 0.64930555555555558
$ws.Range("E52").Value = 0.70833333333333337
$ws.Range("G52").Value = "Validated and verified waveforms for ArithUnit.vhd and Adder.vhd"

# Row 53
$ws.Range("B53").Value = 6977
$ws.Range("C53").Value = 43926
$ws.Range("D53").Value = 0.73958333333333337
$ws.Range("E53").Value = 0.76041666666666663
$ws.Range("G53").Value = "Re-compiled VHDL code for LogicUnit.vhd, Adder.vhd and ArithUnit.vhd"

# Row 54
$ws.Range("B54").Value = 6977
$ws.Range("C54").Value = 43926
$ws.Range("D54").Value = 0.76041666666666663
$ws.Range("E54").Value = 0.77916666666666667
$ws.Range("G54").Value = "Updated all summary files, .vho files and .sdo files."

# Update the view: top-left cell and selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 37
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("G54").Select()
